# Update recomputed TPM-based NATMI ligand-receptor metrics (Nid1-Itgav)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 25.23990433333333
$ws.Range("H2").Value = 75.719713
$ws.Range("I2").Value = 0.05173702626903214
$ws.Range("J2").Value = 0.05173702626903214
$ws.Range("M2").Value = 2.724001666666667
$ws.Range("N2").Value = 8.172005
$ws.Range("O2").Value = 0.04635500474236593
$ws.Range("P2").Value = 0.04635500474236593
$ws.Range("Q2").Value = 68.75354147050722
$ws.Range("R2").Value = 618.7818732345651
$ws.Range("S2").Value = 0.002398270098056896
$ws.Range("T2").Value = 0.002398270098056896

# Row 3
$ws.Range("G3").Value = 25.23990433333333
$ws.Range("H3").Value = 75.719713
$ws.Range("I3").Value = 0.05173702626903214
$ws.Range("J3").Value = 0.05173702626903214
$ws.Range("O3").Value = 0.6912512390256352
$ws.Range("P3").Value = 0.6912512390256351
$ws.Range("Q3").Value = 1025.260831986333
$ws.Range("R3").Value = 9227.347487876998
$ws.Range("S3").Value = 0.0357632835119703
$ws.Range("T3").Value = 0.0357632835119703

# Row 4
$ws.Range("G4").Value = 25.23990433333333
$ws.Range("H4").Value = 75.719713
$ws.Range("I4").Value = 0.05173702626903214
$ws.Range("J4").Value = 0.05173702626903214
$ws.Range("M4").Value = 15.419285
$ws.Range("N4").Value = 46.257855
$ws.Range("O4").Value = 0.2623937562319988
$ws.Range("P4").Value = 0.2623937562319988
$ws.Range("Q4").Value = 389.1812782884016
$ws.Range("R4").Value = 3502.631504595615
$ws.Range("S4").Value = 0.01357547265900494
$ws.Range("T4").Value = 0.01357547265900494

# Row 5
$ws.Range("I5").Value = 0.8454897015965644
$ws.Range("J5").Value = 0.8454897015965646
$ws.Range("M5").Value = 2.724001666666667
$ws.Range("N5").Value = 8.172005
$ws.Range("O5").Value = 0.04635500474236593
$ws.Range("P5").Value = 0.04635500474236593
$ws.Range("Q5").Value = 1123.574651533478
$ws.Range("R5").Value = 10112.1718638013
$ws.Range("S5").Value = 0.0391926791271303
$ws.Range("T5").Value = 0.03919267912713031

# Row 6
$ws.Range("I6").Value = 0.8454897015965644
$ws.Range("J6").Value = 0.8454897015965646
$ws.Range("O6").Value = 0.6912512390256352
$ws.Range("P6").Value = 0.6912512390256351
$ws.Range("S6").Value = 0.5844458038120398
$ws.Range("T6").Value = 0.5844458038120398

# Row 7
$ws.Range("I7").Value = 0.8454897015965644
$ws.Range("J7").Value = 0.8454897015965646
$ws.Range("M7").Value = 15.419285
$ws.Range("N7").Value = 46.257855
$ws.Range("O7").Value = 0.2623937562319988
$ws.Range("P7").Value = 0.2623937562319988
$ws.Range("Q7").Value = 6360.024658858033
$ws.Range("R7").Value = 57240.2219297223
$ws.Range("S7").Value = 0.2218512186573944
$ws.Range("T7").Value = 0.2218512186573944

# Row 8
$ws.Range("G8").Value = 50.137933
$ws.Range("H8").Value = 150.413799
$ws.Range("I8").Value = 0.1027732721344034
$ws.Range("J8").Value = 0.1027732721344034
$ws.Range("M8").Value = 2.724001666666667
$ws.Range("N8").Value = 8.172005
$ws.Range("O8").Value = 0.04635500474236593
$ws.Range("P8").Value = 0.04635500474236593
$ws.Range("Q8").Value = 136.5758130552217
$ws.Range("R8").Value = 1229.182317496995
$ws.Range("S8").Value = 0.004764055517178732
$ws.Range("T8").Value = 0.004764055517178732

# Row 9
$ws.Range("G9").Value = 50.137933
$ws.Range("H9").Value = 150.413799
$ws.Range("I9").Value = 0.1027732721344034
$ws.Range("J9").Value = 0.1027732721344034
$ws.Range("O9").Value = 0.6912512390256352
$ws.Range("P9").Value = 0.6912512390256351
$ws.Range("Q9").Value = 2036.634458783079
$ws.Range("R9").Value = 18329.71012904771
$ws.Range("S9").Value = 0.07104215170162512
$ws.Range("T9").Value = 0.0710421517016251

# Row 10
$ws.Range("G10").Value = 50.137933
$ws.Range("H10").Value = 150.413799
$ws.Range("I10").Value = 0.1027732721344034
$ws.Range("J10").Value = 0.1027732721344034
$ws.Range("M10").Value = 15.419285
$ws.Range("N10").Value = 46.257855
$ws.Range("O10").Value = 0.2623937562319988
$ws.Range("P10").Value = 0.2623937562319988
$ws.Range("Q10").Value = 773.091078237905
$ws.Range("R10").Value = 6957.819704141144
$ws.Range("S10").Value = 0.02696706491559952
$ws.Range("T10").Value = 0.02696706491559952
